# Generate Report for Archive
# - Flip the localization status from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn / de-de status columns) and on each of the
#   per-locale detail sheets ("zh-cn", "de-de") in their Status column.
# - Narrow the (now shorter) status columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511

# --- zh-cn detail sheet: column C is the Status column ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511

# --- de-de detail sheet: column C is the Status column ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
